# issue #5: stock data from json to db
#
# The "股票" (stock) sheet (sheet6.xml) gains three new columns:
#   - "category"     inserted right after "property_category" (new col I)
#   - "source_file"  appended after "legislator_id" (new col M)
#   - "index"        appended after "source_file" (new col N)
#
# This shifts the former I/J/K ("date"/"legislator_name"/"legislator_id")
# columns one place to the right (J/K/L). We also fix a stray bullet
# character in a stock name ("中鋼•" -> "中鋼").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

$lastRow = 15

# --- Header row (row 1): shift I1:K1 -> J1:L1, then fill in new headers ---
$ws.Cells.Item(1, 12).Value = $ws.Cells.Item(1, 11).Value2   # L1 = legislator_id (was K1)
$ws.Cells.Item(1, 11).Copy()
$ws.Cells.Item(1, 12).PasteSpecial(-4122)

$ws.Cells.Item(1, 11).Value = $ws.Cells.Item(1, 10).Value2   # K1 = legislator_name (was J1)

$ws.Cells.Item(1, 10).Value = $ws.Cells.Item(1, 9).Value2    # J1 = date (was I1)

$ws.Cells.Item(1, 9).Value = "category"                      # I1 = category (new)
$ws.Cells.Item(1, 9).Copy()
$ws.Cells.Item(1, 10).PasteSpecial(-4122)
$ws.Cells.Item(1, 11).PasteSpecial(-4122)
$ws.Cells.Item(1, 12).PasteSpecial(-4122)

$ws.Cells.Item(1, 13).Value = "source_file"                  # M1 = source_file (new)
$ws.Cells.Item(1, 9).Copy()
$ws.Cells.Item(1, 13).PasteSpecial(-4122)

$ws.Cells.Item(1, 14).Value = "index"                         # N1 = index (new)
$ws.Cells.Item(1, 9).Copy()
$ws.Cells.Item(1, 14).PasteSpecial(-4122)

# --- Data rows (rows 2-15) ---
for ($r = 2; $r -le $lastRow; $r++) {
    $oldDate = $ws.Cells.Item($r, 9).Value2
    $oldLegislatorName = $ws.Cells.Item($r, 10).Value2
    $oldLegislatorId = $ws.Cells.Item($r, 11).Value2
    $rowIndex = $ws.Cells.Item($r, 1).Value2

    # L = legislator_id (was K)
    $ws.Cells.Item($r, 12).Value = $oldLegislatorId
    $ws.Cells.Item($r, 11).Copy()
    $ws.Cells.Item($r, 12).PasteSpecial(-4122)

    # K = legislator_name (was J)
    $ws.Cells.Item($r, 11).Value = $oldLegislatorName

    # J = date (was I)
    $ws.Cells.Item($r, 10).Value = $oldDate

    # I = category (new)
    $ws.Cells.Item($r, 9).Value = "normal"

    # M = source_file (new)
    $ws.Cells.Item($r, 13).Value = "tmpe6fb1"
    $ws.Cells.Item($r, 9).Copy()
    $ws.Cells.Item($r, 13).PasteSpecial(-4122)

    # N = index (new) -- same as column A for this sheet
    $ws.Cells.Item($r, 14).Value = $rowIndex
    $ws.Cells.Item($r, 9).Copy()
    $ws.Cells.Item($r, 14).PasteSpecial(-4122)
}

# Fix stray bullet character in stock name "中鋼•" -> "中鋼"
for ($r = 2; $r -le $lastRow; $r++) {
    $name = $ws.Cells.Item($r, 2).Value2
    if ($name -eq "中鋼•") {
        $ws.Cells.Item($r, 2).Value = "中鋼"
    }
}

Write-Output "done"
